# RTP - Fluid Assignment (Presentation) edits
#
# This script reproduces, via the PowerPoint COM object model, the
# content edits described by the target diff:
#   1. Slide 9  ("My Approach" title) - split off "Approach" and colour it blue (0070C0)
#   2. Slide 10 ("My Approach" title) - same edit as slide 9
#   3. Slide 2  - "Decided that Muller's people ..." -> "...paper..." (and split into 3 runs)
#   4. Slide 2  - "Less concerned with their attempts to visualise ..." -> "...attempts to
#                  realistically visualise..." (and split into 5 runs)
#   5. Slide 3  - merge the 3-run "SPH provides us ... each particle." sentence into 1 run
#   6. Slide 8  - merge the 3-run "Main stumbling block ... constants!" sentence into 1 run
#   7. Notes of slide 1 - update the cached date/time field text (best effort)
#
# Note: TextRange.Text round-trips "smart" typographic punctuation (curly quotes,
# apostrophes) as their plain ASCII equivalents, so we search using the plain
# ASCII forms but write back using the actual Unicode punctuation so the saved
# OOXML keeps matching the original typography.

$p = $ppt.ActivePresentation

$blue = 12611584  # RGB(0, 112, 192) = 0x0070C0

function Set-TitleApproachColor($slideIndex) {
    $s = $p.Slides.Item($slideIndex)
    $titleShape = $s.Shapes.Item(1)
    $tr = $titleShape.TextFrame.TextRange
    $fullText = $tr.Text
    $marker = "Approach"
    $idx = $fullText.IndexOf($marker)
    if ($idx -ge 0) {
        $sub = $tr.Characters($idx + 1, $marker.Length)
        $sub.Font.Color.RGB = $blue
    }
}

# 1 & 2: "My Approach" titles on slides 9 and 10
Set-TitleApproachColor 9
Set-TitleApproachColor 10

# 3: Slide 2 - "Decided that Muller's people was easier to visualise."
#    -> "Decided that Muller's paper was easier to visualise.", split into 3 runs
$s2 = $p.Slides.Item(2)
$body2 = $s2.Shapes.Item(2).TextFrame.TextRange

$fullText2 = $body2.Text
$oldSentence1 = "Decided that Muller's people was easier to visualise."
$newSentence1 = "Decided that Muller" + [char]0x2019 + "s paper was easier to visualise."
$idx1 = $fullText2.IndexOf($oldSentence1)
$sentRange1 = $body2.Characters($idx1 + 1, $oldSentence1.Length)
$sentRange1.Text = $newSentence1

# Re-split the replaced sentence into 3 runs (no formatting differences, matching
# the target OOXML which shows the same rPr repeated on each run). We force the
# split by re-asserting a font property to its own current value, which causes
# the underlying run to break at that boundary without altering formatting.
$fullText2b = $body2.Text
$idx1b = $fullText2b.IndexOf("Decided that Muller")
$part1a = "Decided that Muller" + [char]0x2019 + "s "
$part1b = "paper was "
$r1a = $body2.Characters($idx1b + 1, $part1a.Length)
$r1a.Font.Underline = $r1a.Font.Underline
$r1b = $body2.Characters($idx1b + 1 + $part1a.Length, $part1b.Length)
$r1b.Font.Underline = $r1b.Font.Underline

# 4: Slide 2 - "Less concerned with their attempts to visualise surface particles."
#    -> "Less concerned with attempts to realistically visualise surface particles.",
#    split into 5 runs
$fullText2c = $body2.Text
$oldSentence2 = "Less concerned with their attempts to visualise surface particles."
$newSentence2 = "Less concerned with attempts to realistically visualise surface particles."
$idx2 = $fullText2c.IndexOf($oldSentence2)
$sentRange2 = $body2.Characters($idx2 + 1, $oldSentence2.Length)
$sentRange2.Text = $newSentence2

$fullText2d = $body2.Text
$idx2b = $fullText2d.IndexOf("Less concerned with")
$p2a = "Less concerned with "
$p2b = "attempts "
$p2c = "to "
$p2d = "realistically visualise "
$r2a = $body2.Characters($idx2b + 1, $p2a.Length)
$r2a.Font.Underline = $r2a.Font.Underline
$r2b = $body2.Characters($idx2b + 1 + $p2a.Length, $p2b.Length)
$r2b.Font.Underline = $r2b.Font.Underline
$r2c = $body2.Characters($idx2b + 1 + $p2a.Length + $p2b.Length, $p2c.Length)
$r2c.Font.Underline = $r2c.Font.Underline
$r2d = $body2.Characters($idx2b + 1 + $p2a.Length + $p2b.Length + $p2c.Length, $p2d.Length)
$r2d.Font.Underline = $r2d.Font.Underline

# 5: Slide 3 - merge "SPH provides us with smoothed approximations of the quantity of
#    fluid "in" each particle." into a single run
$s3 = $p.Slides.Item(3)
$body3 = $s3.Shapes.Item(3).TextFrame.TextRange
$fullText3 = $body3.Text
$oldSentence3 = "SPH provides us with smoothed approximations of the quantity of fluid " + [char]0x22 + "in" + [char]0x22 + " each particle."
$newSentence3 = "SPH provides us with smoothed approximations of the quantity of fluid " + [char]0x201C + "in" + [char]0x201D + " each particle."
$idx3 = $fullText3.IndexOf($oldSentence3)
$sentRange3 = $body3.Characters($idx3 + 1, $oldSentence3.Length)
$sentRange3.Text = $newSentence3

# 6: Slide 8 - merge "Main stumbling block - Fine tuning the fluid constants!" into a single run
$s8 = $p.Slides.Item(8)
$body8 = $s8.Shapes.Item(2).TextFrame.TextRange
$fullText8 = $body8.Text
$sentence8 = "Main stumbling block - Fine tuning the fluid constants!"
$idx8 = $fullText8.IndexOf("Main stumbling block")
$sentRange8 = $body8.Characters($idx8 + 1, $sentence8.Length)
$sentRange8.Text = $sentence8

# 7: Notes of slide 1 - update the cached date/time field text (best effort; the
#    underlying <a:fld> is normally recalculated by the host application rather
#    than being free text, so this may be a no-op under emulation).
$s1 = $p.Slides.Item(1)
$notes1 = $s1.NotesPage
$dateShape = $notes1.Shapes.Item(4)
$dateShape.TextFrame.TextRange.Text = "3/31/2014 1:32 PM"
